$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header rows (Spreadsheet + Test table headers), merged like the other blocks ---
$ws.Range("C64:D64").Merge()
$ws.Range("F64:G64").Merge()
$ws.Range("C64:D64").HorizontalAlignment = -4108
$ws.Range("F64:G64").HorizontalAlignment = -4108

$ws.Range("C64").Value = 'Spreadsheet SpreadsheetResult error5()'
$ws.Range("F64").Value = 'Test error5 error5_test1'

# --- Row 65: column headers for Spreadsheet table (S/F) and Test table ---
$ws.Range("C65").Value = 'S'
$ws.Range("D65").Value = 'F'
$ws.Range("F65").Value = '_description_'
$ws.Range("G65").Value = '_res_.$Result'
$ws.Range("H65").Value = '_error_.code'
$ws.Range("I65").Value = '_error_.msg'
$ws.Range("J65").Value = '_error_.base.code'
$ws.Range("K65").Value = '_error_.base.msg'
$ws.Range("L65").Value = '_error_.nested[0].code'
$ws.Range("M65").Value = '_error_.nested[0].msg'
$ws.Range("N65").Value = '_error_.nested[1].code'
$ws.Range("O65").Value = '_error_.nested[1].msg'

# --- Row 66 ---
$ws.Range("C66").Value = 'base'
$ws.Range("D66").Value = '''=Error(msg="base", code = 17)'
$ws.Range("F66").Value = '_description_'
$ws.Range("G66").Value = '_res_'
$ws.Range("H66").Value = '_error_.code'
$ws.Range("I66").Value = '_error_.msg'
$ws.Range("J66").Value = '_error_.base.code'
$ws.Range("K66").Value = '_error_.base.msg'
$ws.Range("L66").Value = '_error_.nested[0].code'
$ws.Range("M66").Value = '_error_.nested[0].msg'
$ws.Range("N66").Value = '_error_.nested[1].code'
$ws.Range("O66").Value = '_error_.nested[1].msg'

# --- Row 67 ---
$ws.Range("C67").Value = 'x1'
$ws.Range("D67").Value = '''=Error(msg="x1", code = 1)'
$ws.Range("F67").Value = 'fail'

# --- Row 68 ---
$ws.Range("C68").Value = 'x2'
$ws.Range("D68").Value = '''=Error(msg="x2", code = 2)'
$ws.Range("F68").Value = 'fail'
$ws.Range("G68").Value = 'Not achieved'

# --- Row 69 ---
$ws.Range("C69").Value = 'nested'
$ws.Range("D69").Value = '''=add($x1, $x2)'
$ws.Range("F69").Value = 'pass'
$ws.Range("H69").Value = 42
$ws.Range("I69").Value = 'total'
$ws.Range("J69").Value = 17
$ws.Range("K69").Value = 'base'
$ws.Range("L69").Value = 1
$ws.Range("M69").Value = 'x1'
$ws.Range("N69").Value = 2
$ws.Range("O69").Value = 'x2'

# --- Row 70 ---
$ws.Range("C70").Value = 'error'
$ws.Range("D70").Value = '''=new Error("total", 42, $nested, $base)'
$ws.Range("F70").Value = 'fail'
$ws.Range("I70").Value = 'total'
$ws.Range("J70").Value = 17
$ws.Range("K70").Value = 'base'
$ws.Range("L70").Value = 1
$ws.Range("M70").Value = 'x1'
$ws.Range("N70").Value = 2
$ws.Range("O70").Value = 'x2'

# --- Row 71 ---
$ws.Range("C71").Value = 'Step'
$ws.Range("D71").Value = '''= error($error)'
$ws.Range("F71").Value = 'fail'
$ws.Range("H71").Value = 42
$ws.Range("J71").Value = 17
$ws.Range("K71").Value = 'base'
$ws.Range("L71").Value = 1
$ws.Range("M71").Value = 'x1'
$ws.Range("N71").Value = 2
$ws.Range("O71").Value = 'x2'

# --- Row 72 ---
$ws.Range("C72").Value = 'Result'
$ws.Range("D72").Value = 'Not achieved'
$ws.Range("F72").Value = 'fail'
$ws.Range("H72").Value = 42
$ws.Range("I72").Value = 'total'
$ws.Range("K72").Value = 'base'
$ws.Range("L72").Value = 1
$ws.Range("M72").Value = 'x1'
$ws.Range("N72").Value = 2
$ws.Range("O72").Value = 'x2'

# --- Row 73 ---
$ws.Range("F73").Value = 'fail'
$ws.Range("H73").Value = 42
$ws.Range("I73").Value = 'total'
$ws.Range("J73").Value = 17
$ws.Range("L73").Value = 1
$ws.Range("M73").Value = 'x1'
$ws.Range("N73").Value = 2
$ws.Range("O73").Value = 'x2'

# --- Row 74 ---
$ws.Range("F74").Value = 'fail'
$ws.Range("H74").Value = 42
$ws.Range("I74").Value = 'total'
$ws.Range("J74").Value = 17
$ws.Range("K74").Value = 'base'
$ws.Range("M74").Value = 'x1'
$ws.Range("N74").Value = 2
$ws.Range("O74").Value = 'x2'

# --- Row 75 ---
$ws.Range("F75").Value = 'fail'
$ws.Range("H75").Value = 42
$ws.Range("I75").Value = 'total'
$ws.Range("J75").Value = 17
$ws.Range("K75").Value = 'base'
$ws.Range("L75").Value = 1
$ws.Range("N75").Value = 2
$ws.Range("O75").Value = 'x2'

# --- Row 76 ---
$ws.Range("C76").Value = 'Datatype Error'
$ws.Range("F76").Value = 'fail'
$ws.Range("H76").Value = 42
$ws.Range("I76").Value = 'total'
$ws.Range("J76").Value = 17
$ws.Range("K76").Value = 'base'
$ws.Range("L76").Value = 1
$ws.Range("M76").Value = 'x1'
$ws.Range("O76").Value = 'x2'

# --- Row 77 ---
$ws.Range("C77").Value = 'String'
$ws.Range("D77").Value = 'msg'
$ws.Range("F77").Value = 'fail'
$ws.Range("H77").Value = 42
$ws.Range("I77").Value = 'total'
$ws.Range("J77").Value = 17
$ws.Range("K77").Value = 'base'
$ws.Range("L77").Value = 1
$ws.Range("M77").Value = 'x1'
$ws.Range("N77").Value = 2

# --- Row 78 ---
$ws.Range("C78").Value = 'Integer'
$ws.Range("D78").Value = 'code'
$ws.Range("F78").Value = 'fail'
$ws.Range("H78").Value = 42
$ws.Range("I78").Value = 'total'
$ws.Range("J78").Value = 17
$ws.Range("K78").Value = 'base'
$ws.Range("L78").Value = 1
$ws.Range("M78").Value = 'x1'

# --- Row 79 ---
$ws.Range("C79").Value = 'Error[]'
$ws.Range("D79").Value = 'nested'
$ws.Range("F79").Value = 'fail'
$ws.Range("H79").Value = 42
$ws.Range("I79").Value = 'total'
$ws.Range("J79").Value = 17
$ws.Range("K79").Value = 'base'

# --- Row 80 ---
$ws.Range("C80").Value = 'Error'
$ws.Range("D80").Value = 'base'
$ws.Range("F80").Value = 'fail'
$ws.Range("H80").Value = 42
$ws.Range("I80").Value = 'total'

# --- Row 81 ---
$ws.Range("F81").Value = 'fail'
$ws.Range("H81").Value = 42

# --- Row 82 ---
$ws.Range("F82").Value = 'fail'
$ws.Range("I82").Value = 'total'

# --- Column width tweaks (match the widened columns from the authoring tool) ---
$ws.Range("D1").ColumnWidth = 34.13
$ws.Range("F1").ColumnWidth = 21.24
$ws.Range("H1").ColumnWidth = 12.44
$ws.Range("I1").ColumnWidth = 12.44
$ws.Range("J1").ColumnWidth = 16.6
$ws.Range("K1").ColumnWidth = 16.13
$ws.Range("S1").ColumnWidth = 12.03

# --- View state: scroll to show the new table, select F82 like the source diff ---
$ws.Application.ActiveWindow.ScrollRow = 28
$ws.Range("F82").Select()
